# Refresh cryptos list: update Price (D) / Volume(1h) (E) columns for rows 2-31
# and 41-51, and re-splice rows 32-40 (Binance-PegBSC-USD moved up to rank 30,
# the rest shifting down one row) with the data from the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values that look like plain numbers must be pre-formatted as Text so
# Excel keeps the exact original string (e.g. trailing zeros) instead of coercing
# them to a number.

$ws.Range("D2").Value = '90.997.47'
$ws.Range("E2").Value = '  -0.36%  '

$ws.Range("D3").Value = '3.148.44'
$ws.Range("E3").Value = '  +1.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.47'
$ws.Range("E5").Value = '  +8.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '643.06'
$ws.Range("E6").Value = '  +3.20%  '

$ws.Range("E7").Value = '  +10.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.368'
$ws.Range("E8").Value = '  -2.65%  '

$ws.Range("E9").Value = '  +0.07%  '

$ws.Range("D10").Value = '3.133.91'
$ws.Range("E10").Value = '  +0.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.724'
$ws.Range("E11").Value = '  +0.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.197'
$ws.Range("E12").Value = '  +3.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.66'
$ws.Range("E13").Value = '  +5.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000252'
$ws.Range("E14").Value = '  -1.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.60'
$ws.Range("E15").Value = '  +3.59%  '

$ws.Range("D16").Value = '90.687.65'
$ws.Range("E16").Value = '  -0.38%  '

$ws.Range("D17").Value = '3.731.36'
$ws.Range("E17").Value = '  +1.07%  '

$ws.Range("D18").Value = '3.148.33'
$ws.Range("E18").Value = '  +0.97%  '

$ws.Range("E19").Value = '  -0.31%  '

$ws.Range("E20").Value = '  -0.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.47'
$ws.Range("E21").Value = '  +2.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '450.31'
$ws.Range("E22").Value = '  +3.24%  '

$ws.Range("E23").Value = '  +9.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.04'
$ws.Range("E24").Value = '  +2.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.02'
$ws.Range("E25").Value = '  -2.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.23'
$ws.Range("E26").Value = '  +5.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.44'
$ws.Range("E27").Value = '  +1.38%  '

$ws.Range("E28").Value = '  +0.62%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("E30").Value = '  +7.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.160'
$ws.Range("E31").Value = '  -4.54%  '

$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.966'
$ws.Range("E32").Value = '  +7.89%  '

$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.204'
$ws.Range("E33").Value = '  +33.30%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.10'
$ws.Range("E34").Value = '  +14.72%  '

$ws.Range("B35").Value = 'dogwifhat'
$ws.Range("C35").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.89'
$ws.Range("E35").Value = '  +3.18%  '

$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '518.16'
$ws.Range("E36").Value = '  -1.30%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.151'
$ws.Range("E37").Value = '  +3.13%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.19'
$ws.Range("E38").Value = '  +0.35%  '

$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.95'
$ws.Range("E39").Value = '  +4.64%  '

$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.31'
$ws.Range("E40").Value = '  +0.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.422'
$ws.Range("E41").Value = '  +5.42%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.20'
$ws.Range("E42").Value = '  -0.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0864'
$ws.Range("E43").Value = '  -5.03%  '

$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.40'
$ws.Range("E45").Value = '  +47.00%  '

$ws.Range("E46").Value = '  +1.07%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.706'
$ws.Range("E47").Value = '  +12.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '151.39'
$ws.Range("E48").Value = '  +1.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '45.69'
$ws.Range("E49").Value = '  +3.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.58'
$ws.Range("E50").Value = '  +7.96%  '

$ws.Range("E51").Value = '  +4.02%  '
